# Auto update stock data
# Updates the "as of" date (column A) and EBITDA figure (column B) for each
# company block in the risk-scores sheet. Cells are stored as text, so the
# target range is formatted as Text before the new values are written -
# this stops Excel from re-interpreting "2025/12/31" as a date serial or
# "6.77" as a floating point number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new EBITDA value (row 38 keeps its original EBITDA figure)
$updates = @{
    2  = "6.77"
    8  = "8.57"
    14 = "3.02"
    20 = "12.79"
    26 = "11.33"
    32 = "27.46"
    38 = $null
    44 = "11.11"
    50 = "11.29"
    56 = "31.40"
    62 = "11.55"
    68 = "12.86"
    74 = "16.58"
}

foreach ($row in $updates.Keys) {
    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2025/12/31"

    $newB = $updates[$row]
    if ($newB -ne $null) {
        $ebitdaCell = $ws.Cells.Item($row, 2)
        $ebitdaCell.NumberFormat = "@"
        $ebitdaCell.Value = $newB
    }
}
